$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New tweet data: replaces the previous content that lived in row 2, and is
# also duplicated into a brand-new row 3 (same six values, same shared
# strings once written).
$rowValues  = @(
    "['Naty está encantada', '@fadadesaturno']",
    "30",
    "77",
    "1.688",
    "28 mil",
    "24 de ago Meu pai sempre apoiou meu sonho de ser escritora, e hoje tirei uma foto LINDA dele com meu livro!!! 30"
)
# Columns whose text would otherwise be auto-coerced into a number by Excel
# (e.g. "30" -> 30). Force them to stay text, then drop the format change so
# no extra cell style sticks around.
$forceText = @(2, 3, 4)

foreach ($row in 2, 3) {
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($forceText -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $rowValues[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $rowValues[$i]
        }
    }
}
